{"js": "// Apply the documented text edits to the body using Word JS API search/replace.\nconst body = context.document.body;\n\nconst replacements = [\n  { find: \"nattsk\u00e4rra\", replace: \"nattsk\u00e4rra och spillkr\u00e5ka\", all: true },\n  { find: \"Detta \u00e4r en prioriterad art\", replace: \"Dessa \u00e4r prioriterade arter\", all: false },\n  { find: \"denna art\", replace: \"dessa arter\", all: false },\n  { find: \"arten\", replace: \"arterna\", all: false },\n  { find: \"2026-02-24\", replace: \"2026-02-25\", all: false },\n];\n\nfor (const { find, replace } of replacements) {\n  const found = body.search(find, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the documented text edits via Word's Find/Replace (COM object model).\n$d = $word.ActiveDocument\n\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\nfunction Replace-Text($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Execute($findText, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $replaceText, $wdReplaceAll)\n}\n\nReplace-Text \"nattsk\u00e4rra\" \"nattsk\u00e4rra och spillkr\u00e5ka\"\nReplace-Text \"Detta \u00e4r en prioriterad art\" \"Dessa \u00e4r prioriterade arter\"\nReplace-Text \"denna art\" \"dessa arter\"\nReplace-Text \"arten\" \"arterna\"\nReplace-Text \"2026-02-24\" \"2026-02-25\"\n"}
